$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "dgfhdfdfdfdfdfdfdf"
$ws.Range("C2").Value = "Sajad@gmail.com"
